# "Correct the weird bug": the timesheet export was missing a
# "Date of Last Update" column. Insert two new columns between
# "Time Worked (Minutes)" (H) and "Defects" (old I, now shifted to K):
#
#   I -> "Date of Last Update" (re-uses the date that used to sit in the
#        old "Defects" column, i.e. the second timestamp)
#   J -> a duplicate "TicketID" column (the actual bug being preserved/
#        reproduced - it mirrors column A's header & value)
#
# The old "Defects" / "Timesheet - Related" columns shift right to K/L,
# and the (until-now-empty) "Defects" data cell is filled in with
# "Needed / Missing".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert two blank columns at I:J (old I/J -> K/L).
$ws.Columns("I:J").Insert()

# --- Header row ---
$ws.Range("I1").Value = "Date of Last Update"
$ws.Range("J1").Value = "TicketID"
$ws.Range("K1").Value = "Defects"

# --- Data row ---
$ws.Range("I2").Value = "2017-08-15T09:31:24.763000"

# Copy column A's value into J2 so it stays text ("252"), matching the
# duplicated-TicketID bug exactly instead of being re-interpreted as a number.
$ws.Range("A2").Copy($ws.Range("J2"))

$ws.Range("K2").Value = "Needed / Missing"
